# Bump the "Förändrad" (Changed) date in column C for all data rows
# (rows 2-34) from 46061 (2026-02-08) to 46062 (2026-02-09).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 34; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46061) {
        $cell.Value = 46062
    }
}
